$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "contenttable" (sheet6.xml): add a new "taskFile" column and
# a new row describing a file-type task (#43 "Fixed API for task file").
# -----------------------------------------------------------------
$content = $wb.Worksheets.Item("contenttable")

# New header cell T1 = "taskFile" - copy formatting from the existing
# bold header style used by S1 ("questions") so no stray new font/style
# is created.
$content.Range("T1").Value = "taskFile"
$content.Range("S1").Copy()
$content.Range("T1").PasteSpecial(-4122)

# New data row (row 8) describing the sample file task.
$content.Range("A8").Value = 7
$content.Range("B8").Value = "fileTask"
$content.Range("C8").Value = "teacher"
$content.Range("D8").Value = "file test"
$content.Range("E8").Value = "file"
$content.Range("F8").Value = 60
$content.Range("G8").Value = 0
$content.Range("I8").Value = 2
$content.Range("J8").Value = 1
$content.Range("K8").Value = 2
$content.Range("L8").Value = "2020-12-12 00:00:00"
$content.Range("M8").Value = 0
$content.Range("N8").Value = "2020-12-12 00:00:00"
$content.Range("O8").Value = "2020-12-12 00:00:00"
$content.Range("P8").Value = 2
$content.Range("Q8").Value = 7
$content.Range("T8").Value = ".\files\task\admin_c_1651829254.4129157\6磁场-21_admin_tsk_1651244631_admin_tsk_1651829254.4139159.pdf"

# -----------------------------------------------------------------
# Sheet "usercontenttable" (sheet11.xml): add the matching join record
# for the newcomer who was assigned the new fileTask (content id 7).
# -----------------------------------------------------------------
$userContent = $wb.Worksheets.Item("usercontenttable")

$userContent.Range("B28").Value = "newcomer"
$userContent.Range("C28").Value = "7"
$userContent.Range("D28").Value = "0"
$userContent.Range("E28").Value = "2022-5-7 00:00:05"
$userContent.Range("F28").Value = "2022-5-7 00:00:05"
$userContent.Range("G28").Value = "2022-5-7 00:00:05"
$userContent.Range("H28").Value = "teacher"
$userContent.Range("I28").Value = "1"
$userContent.Range("J28").Value = "0"
$userContent.Range("K28").Value = "28"
$userContent.Range("L28").Value = "27"
